$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("n9",  "n9_IMG_3180.jpeg",  "'True", "no_meltpatch", "negative"),
    @("n10", "n10_IMG_3177.jpeg", "'True", "no_meltpatch", "negative"),
    @("n11", "n11_IMG_3175.jpeg", "'True", "no_meltpatch", "negative"),
    @("n12", "n12_IMG_3176.jpeg", "'True", "no_meltpatch", "negative")
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $cell = $ws.Cells.Item($rowNum, $col)
        $cell.Value = $rowData[$col - 1]
        $cell.Style = "Normal"
    }
}

$wb.Save()
